# Apply the "LessonTree" content update: fill in the lesson step details
# that were missing (row 2's starting prompt/response, row 3's assistant
# reply, and row 11's example code), then leave the selection on the
# newly-edited cell B11, matching the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "I know a little bit of the for loop"
$ws.Range("B2").Value = "Tell the user it is great they know about for loops and ask the user to write the syntax of the for loop"
$ws.Range("B3").Value = "Tell the user you are happy they know how to write a for loop. Ask the user to write a  for loop that will display numbers 1 to 10."
$ws.Range("A11").Value = "int numbers[]=[45,67,77,56,78];`nfor(i=0;i<5;i++)`n{printf(`"%d`",i)}"

# Move the selection/scroll position to B11, where the new content was added.
$ws.Range("A10").Select()
$ws.Range("B11").Select()
